$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.159.55'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '1.669.48'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5235'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.34%  '

$ws.Range("E7").Value = '  -0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2613'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.98%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06354'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.12'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07540'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.22%  '

$ws.Range("D12").Value = '1.671.60'
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.436'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5442'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.62%  '

$ws.Range("D15").Value = '0.0₅8021'
$ws.Range("E15").Value = '  -1.70%  '

$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").Value = '26.196.00'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.746'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.262'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.66'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1232'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.453'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.77'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06264'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.363'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.60%  '

$ws.Range("E30").Value = '  -1.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.496'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.419'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.645'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.002'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.04%  '

$ws.Range("B35").Value = 'MXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.760'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.32%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.392'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5958'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.66%  '

$ws.Range("D38").Value = '1.111.61'
$ws.Range("E38").Value = '  +0.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.065'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01606'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8585'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.73%  '

$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.79'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.63%  '

$ws.Range("D44").Value = '1.816.08'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("D45").Value = '0.0₈110'
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.58'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.069'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05243'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4236'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.924'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.09%  '
